$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (row, dateSerial, B, C, D)
$data = @(
    @(245, 44319, 0, 5, 240.3846153846154),
    @(246, 44320, 0, 3, 144.2307692307692),
    @(247, 44321, 0, 3, 144.2307692307692)
)

foreach ($row in $data) {
    $r = $row[0]

    # Copy formatting from the last existing data row (244) for column A
    # so the new date cell reuses the existing date style (s="2").
    $ws.Cells.Item(244, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

$excel.CutCopyMode = $false
